$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new customer record as row 18 (mirrors the layout of the
# existing rows: columns A,B,C,E,F,G,H populated, column D left empty).
#
# The values are entered as string-literal formulas first (so Excel keeps
# them as text -- preserving e.g. the leading zero in the CPF number --
# without flipping the workbook into "number" type) and are then converted
# back to plain static values via copy / paste-special-values, leaving no
# formula behind.
$ws.Cells.Item(18, 1).Formula = '="BRUNO DE FRAGA"'
$ws.Cells.Item(18, 2).Formula = '="123123123"'
$ws.Cells.Item(18, 3).Formula = '="02370945095"'
$ws.Cells.Item(18, 5).Formula = '="92320-195"'
$ws.Cells.Item(18, 6).Formula = '="joanues@gmail.com"'
$ws.Cells.Item(18, 7).Formula = '="5154548"'
$ws.Cells.Item(18, 8).Formula = '="Rua 3 Pinheiros I, 27"'

$newRow = $ws.Range("A18:H18")
$newRow.Copy()
$newRow.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false
